$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date), Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado and Precio $/Kg values between row 2 and row 4.

$ws.Range("D2").Value = 44280
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 500

$ws.Range("D4").Value = 44277
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 550
